{"js": "// Update the date line and the twenty-five division problems in the table\n// to the new values from the commit, via search+replace on exact text.\nconst replacements = [\n  [\"2026-02-06 Friday\", \"2026-02-07 Saturday\"],\n  [\"792\u00f78=\", \"749\u00f76=\"],\n  [\"223\u00f73=\", \"843\u00f75=\"],\n  [\"511\u00f79=\", \"463\u00f73=\"],\n  [\"567\u00f77=\", \"373\u00f76=\"],\n  [\"451\u00f72=\", \"667\u00f72=\"],\n  [\"534\u00f77=\", \"896\u00f72=\"],\n  [\"562\u00f78=\", \"305\u00f73=\"],\n  [\"347\u00f77=\", \"169\u00f75=\"],\n  [\"356\u00f79=\", \"753\u00f76=\"],\n  [\"354\u00f75=\", \"830\u00f73=\"],\n  [\"973\u00f77=\", \"419\u00f75=\"],\n  [\"317\u00f75=\", \"544\u00f77=\"],\n  [\"264\u00f72=\", \"457\u00f72=\"],\n  [\"756\u00f73=\", \"299\u00f77=\"],\n  [\"292\u00f77=\", \"544\u00f76=\"],\n  [\"663\u00f78=\", \"119\u00f79=\"],\n  [\"347\u00f74=\", \"988\u00f79=\"],\n  [\"584\u00f73=\", \"462\u00f75=\"],\n  [\"309\u00f74=\", \"183\u00f79=\"],\n  [\"663\u00f76=\", \"961\u00f77=\"],\n  [\"431\u00f74=\", \"494\u00f74=\"],\n  [\"260\u00f76=\", \"567\u00f72=\"],\n  [\"167\u00f73=\", \"783\u00f75=\"],\n  [\"833\u00f76=\", \"102\u00f78=\"],\n  [\"765\u00f76=\", \"230\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the twenty-five division problems in the table\n# to the new values from the commit, via Find/Replace on exact text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"2026-02-06 Friday\"; New=\"2026-02-07 Saturday\"},\n    @{Old=\"792\u00f78=\"; New=\"749\u00f76=\"},\n    @{Old=\"223\u00f73=\"; New=\"843\u00f75=\"},\n    @{Old=\"511\u00f79=\"; New=\"463\u00f73=\"},\n    @{Old=\"567\u00f77=\"; New=\"373\u00f76=\"},\n    @{Old=\"451\u00f72=\"; New=\"667\u00f72=\"},\n    @{Old=\"534\u00f77=\"; New=\"896\u00f72=\"},\n    @{Old=\"562\u00f78=\"; New=\"305\u00f73=\"},\n    @{Old=\"347\u00f77=\"; New=\"169\u00f75=\"},\n    @{Old=\"356\u00f79=\"; New=\"753\u00f76=\"},\n    @{Old=\"354\u00f75=\"; New=\"830\u00f73=\"},\n    @{Old=\"973\u00f77=\"; New=\"419\u00f75=\"},\n    @{Old=\"317\u00f75=\"; New=\"544\u00f77=\"},\n    @{Old=\"264\u00f72=\"; New=\"457\u00f72=\"},\n    @{Old=\"756\u00f73=\"; New=\"299\u00f77=\"},\n    @{Old=\"292\u00f77=\"; New=\"544\u00f76=\"},\n    @{Old=\"663\u00f78=\"; New=\"119\u00f79=\"},\n    @{Old=\"347\u00f74=\"; New=\"988\u00f79=\"},\n    @{Old=\"584\u00f73=\"; New=\"462\u00f75=\"},\n    @{Old=\"309\u00f74=\"; New=\"183\u00f79=\"},\n    @{Old=\"663\u00f76=\"; New=\"961\u00f77=\"},\n    @{Old=\"431\u00f74=\"; New=\"494\u00f74=\"},\n    @{Old=\"260\u00f76=\"; New=\"567\u00f72=\"},\n    @{Old=\"167\u00f73=\"; New=\"783\u00f75=\"},\n    @{Old=\"833\u00f76=\"; New=\"102\u00f78=\"},\n    @{Old=\"765\u00f76=\"; New=\"230\u00f77=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
